# Weekly update: insert a new price record for "Jengibre" (Terminal La
# Palmera de La Serena) as row 119, pushing the existing rows 119-134
# down to 120-135 (a new most-recent week's data point is added at the
# top of the date-descending history, just like previous updates did).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 119:134 down to 120:135, leaving a blank row 119 behind.
$ws.Rows("119:119").Insert()

# Populate the new row 119 with this week's record.
$ws.Range("A119").Value = 8
$ws.Range("B119").Value = "Terminal La Palmera de La Serena"
$ws.Range("C119").Value = "Coquimbo"
$ws.Range("D119").Value = 45077
$ws.Range("E119").Value = 4
$ws.Range("F119").Value = 100114007
$ws.Range("G119").Value = "Jengibre"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 360
$ws.Range("K119").Value = 17000
$ws.Range("L119").Value = 18000
$ws.Range("M119").Value = 17500
$ws.Range("N119").Value = "$/caja 13 kilos"
$ws.Range("O119").Value = "Perú"
$ws.Range("P119").Value = 1346
$ws.Range("Q119").Value = 13
$ws.Range("R119").Value = "Hortaliza"

# Match the date-cell number format used by the rest of column D.
$ws.Range("D119").NumberFormat = $ws.Range("D120").NumberFormat
